# Fixed issue #13 Permitir que en los ficheros de metadatos dos columnas se
# puedan relacionar para crear SKOS jerárquicos.
#
# A new row is inserted right after the header row (row 1). It holds the
# "machine name" (slug) of each header in row 1: lowercase, accents
# stripped, spaces/parentheses turned into hyphens. Every row that used to
# sit at row 2+ shifts down by one.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new, blank row at position 2 - pushes existing rows 2..5 down to 3..6.
$ws.Rows.Item(2).Insert()

# The old row 5 only held "mapping-ano.xlsx" (J5); it is now at row 6 and is
# no longer needed, so drop it entirely (rows stay 1..5, same as before).
$ws.Rows.Item(6).Delete()

$slugs = @(
    "mes-codigo",
    "ccaa-nombre",
    "comarca-nombre",
    "mes-nombre",
    "grupo-de-tipo-de-jornada",
    "comarca-codigo",
    "numero-de-contratos",
    "provincia-codigo",
    "provincia-nombre",
    "ano",
    "sexo",
    "mes-y-ano"
)

for ($i = 0; $i -lt $slugs.Length; $i++) {
    $col = $i + 1
    $ws.Cells.Item(2, $col).Value = $slugs[$i]
}
